# fix: improve debug display for time_blocks in session_state
#
# 1) Add a new "Meta" worksheet (key/value store) after Sheet1 with the
#    time_blocks debug info.
# 2) Extend Sheet1 with new tracking columns (O:U) used by the reminder /
#    status-tracking feature, filling existing rows with blank defaults.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Extend Sheet1 with the new O:U columns -----------------------------

# Clone the existing header style (bold + border, same as N1) onto the new
# header cells so they stay on style index 1 instead of minting a new one.
$ws1.Range("N1").Copy()
$ws1.Range("O1:U1").PasteSpecial(-4122)

$ws1.Range("O1").Value = "Patient ID"
$ws1.Range("P1").Value = "REMINDER_SNOOZE_UNTIL"
$ws1.Range("Q1").Value = "REMINDER_DISMISSED"
$ws1.Range("R1").Value = "STATUS_CHANGED_AT"
$ws1.Range("S1").Value = "ACTUAL_START_AT"
$ws1.Range("T1").Value = "ACTUAL_END_AT"
$ws1.Range("U1").Value = "STATUS_LOG"

# Fill the data rows (2-6): every new column starts blank except
# REMINDER_DISMISSED (Q), which defaults to FALSE.
for ($row = 2; $row -le 6; $row++) {
    $ws1.Cells.Item($row, 17).Value = $false   # Q: REMINDER_DISMISSED
}

# --- 2. Add the "Meta" key/value worksheet ----------------------------------

$wsMeta = $wb.Worksheets.Add($null, $ws1)
$wsMeta.Name = "Meta"

# Reuse the same bold header style (index 1) as Sheet1's header row.
$ws1.Range("N1").Copy()
$wsMeta.Range("A1:B1").PasteSpecial(-4122)

$wsMeta.Range("A1").Value = "key"
$wsMeta.Range("B1").Value = "value"

$wsMeta.Range("A2").Value = "time_blocks"
$wsMeta.Range("B2").Value = '[{"assistant": "BABU", "date": "2025-12-27", "reason": "Backend Work", "start_time": "20:46", "end_time": "20:50"}]'

$wsMeta.Range("A3").Value = "time_blocks_updated_at"
$wsMeta.Range("B3").Value = "2025-12-27T02:24:24.675570+05:30"

# Restore the original active sheet/selection (adding Meta made it active).
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
